# New crime data collected — refresh the weekly CompStat report
# (Volume/Number header, reporting week dates, and the Week-to-Date /
# 28-Day / Year-to-Date / 2-Year / 14-Year / 31-Year crime figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header: Volume / Number and reporting week dates ----
$ws.Range("A8").Value2 = "Volume 31   Number  20"
$ws.Range("C9").Value2 = "Report Covering the Week  5/13/2024  Through  5/19/2024"

# ---- Row 15: Rape ----
# C15 switches from the "no data" placeholder to an actual count, so make
# sure it carries the same numeric format used elsewhere in the column.
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value2 = 1
$ws.Range("D15").Value2 = 2
$ws.Range("E15").Value2 = -50
$ws.Range("G15").Value2 = 3
$ws.Range("H15").Value2 = -66.666666666666
$ws.Range("I15").Value2 = 13
$ws.Range("J15").Value2 = 5
$ws.Range("K15").Value2 = 160
$ws.Range("L15").Value2 = 44.444444444444
$ws.Range("M15").Value2 = 550
$ws.Range("N15").Value2 = 62.5

# ---- Row 16: Robbery ----
$ws.Range("C16").Value2 = 5
$ws.Range("D16").Value2 = 10
$ws.Range("E16").Value2 = -50
$ws.Range("F16").Value2 = 29
$ws.Range("G16").Value2 = 37
$ws.Range("H16").Value2 = -21.621621621621
$ws.Range("I16").Value2 = 138
$ws.Range("J16").Value2 = 194
$ws.Range("K16").Value2 = -28.865979381443
$ws.Range("L16").Value2 = -32.352941176470
$ws.Range("M16").Value2 = 146.428571428571
$ws.Range("N16").Value2 = -85.639958376690

# ---- Row 17: Fel. Assault ----
$ws.Range("C17").Value2 = 7
$ws.Range("D17").Value2 = 4
$ws.Range("E17").Value2 = 75
$ws.Range("G17").Value2 = 26
$ws.Range("H17").Value2 = 61.538461538461
$ws.Range("I17").Value2 = 210
$ws.Range("J17").Value2 = 198
$ws.Range("K17").Value2 = 6.060606060606
$ws.Range("L17").Value2 = 56.716417910447
$ws.Range("M17").Value2 = 176.315789473684
$ws.Range("N17").Value2 = -20.152091254752

# ---- Row 18: Burglary ----
$ws.Range("C18").Value2 = 5
$ws.Range("D18").Value2 = 5
$ws.Range("E18").Value2 = 0
$ws.Range("F18").Value2 = 34
$ws.Range("G18").Value2 = 38
$ws.Range("H18").Value2 = -10.526315789473
$ws.Range("I18").Value2 = 151
$ws.Range("J18").Value2 = 180
$ws.Range("K18").Value2 = -16.111111111111
$ws.Range("L18").Value2 = -41.698841698841
$ws.Range("M18").Value2 = 20.8
$ws.Range("N18").Value2 = -85.108481262327

# ---- Row 19: Gr. Larceny ----
$ws.Range("C19").Value2 = 44
$ws.Range("D19").Value2 = 54
$ws.Range("E19").Value2 = -18.518518518518
$ws.Range("F19").Value2 = 159
$ws.Range("G19").Value2 = 214
$ws.Range("H19").Value2 = -25.700934579439
$ws.Range("I19").Value2 = 798
$ws.Range("J19").Value2 = 947
$ws.Range("K19").Value2 = -15.733896515311
$ws.Range("L19").Value2 = 4.177545691906
$ws.Range("M19").Value2 = -7.101280558789
$ws.Range("N19").Value2 = -78.046767537826

# ---- Row 20: G.L.A. ----
$ws.Range("C20").Value2 = 1
$ws.Range("D20").Value2 = 1
$ws.Range("I20").Value2 = 20
$ws.Range("J20").Value2 = 25
$ws.Range("K20").Value2 = -20
$ws.Range("L20").Value2 = -31.034482758620
$ws.Range("M20").Value2 = 122.222222222222
$ws.Range("N20").Value2 = -87.096774193548

# ---- Row 21: TOTAL ----
$ws.Range("C21").Value2 = 63
$ws.Range("D21").Value2 = 76
$ws.Range("E21").Value2 = -17.105263157894
$ws.Range("F21").Value2 = 271
$ws.Range("G21").Value2 = 327
$ws.Range("H21").Value2 = -17.125382262996
$ws.Range("I21").Value2 = 1332
$ws.Range("J21").Value2 = 1550
$ws.Range("K21").Value2 = -14.064516129032
$ws.Range("L21").Value2 = -5.263157894736
$ws.Range("M21").Value2 = 18.189884649512
$ws.Range("N21").Value2 = -77.954319761668

# ---- Row 22: Transit ----
# C22 also switches from the placeholder to a real count.
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value2 = 3
$ws.Range("D22").Value2 = 10
$ws.Range("E22").Value2 = -70
$ws.Range("F22").Value2 = 12
$ws.Range("G22").Value2 = 26
$ws.Range("H22").Value2 = -53.846153846153
$ws.Range("I22").Value2 = 66
$ws.Range("J22").Value2 = 96
$ws.Range("K22").Value2 = -31.25
$ws.Range("L22").Value2 = -12
$ws.Range("M22").Value2 = 32

# ---- Row 24: Petit Larceny ----
$ws.Range("C24").Value2 = 98
$ws.Range("D24").Value2 = 73
$ws.Range("E24").Value2 = 34.246575342465
$ws.Range("F24").Value2 = 348
$ws.Range("G24").Value2 = 323
$ws.Range("H24").Value2 = 7.739938080495
$ws.Range("I24").Value2 = 1604
$ws.Range("J24").Value2 = 1412
$ws.Range("K24").Value2 = 13.597733711048
$ws.Range("L24").Value2 = 46.216955332725
$ws.Range("M24").Value2 = -9.070294784580

# ---- Row 25: Retail Theft ----
$ws.Range("C25").Value2 = 89
$ws.Range("D25").Value2 = 72
$ws.Range("E25").Value2 = 23.611111111111
$ws.Range("F25").Value2 = 302
$ws.Range("G25").Value2 = 316
$ws.Range("H25").Value2 = -4.430379746835
$ws.Range("I25").Value2 = 1411
$ws.Range("J25").Value2 = 1331
$ws.Range("K25").Value2 = 6.010518407212
$ws.Range("L25").Value2 = 32.488262910798

# ---- Row 26: Misd. Assault ----
$ws.Range("C26").Value2 = 24
$ws.Range("D26").Value2 = 22
$ws.Range("E26").Value2 = 9.090909090909
$ws.Range("F26").Value2 = 77
$ws.Range("G26").Value2 = 86
$ws.Range("H26").Value2 = -10.465116279069
$ws.Range("I26").Value2 = 383
$ws.Range("J26").Value2 = 379
$ws.Range("K26").Value2 = 1.055408970976
$ws.Range("L26").Value2 = 20.440251572327
$ws.Range("M26").Value2 = 61.603375527426

# ---- Row 27: UCR Rape* ----
# C27 also switches from the placeholder to a real count.
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value2 = 1
$ws.Range("E27").Value2 = -50
$ws.Range("F27").Value2 = 2
$ws.Range("G27").Value2 = 4
$ws.Range("I27").Value2 = 17
$ws.Range("J27").Value2 = 10
$ws.Range("K27").Value2 = 70
$ws.Range("L27").Value2 = 70

# ---- Row 28: Other Sex Crimes ----
$ws.Range("D28").Value2 = 4
$ws.Range("E28").Value2 = -25
$ws.Range("F28").Value2 = 13
$ws.Range("G28").Value2 = 13
$ws.Range("H28").Value2 = 0
$ws.Range("I28").Value2 = 77
$ws.Range("J28").Value2 = 82
$ws.Range("K28").Value2 = -6.097560975609
$ws.Range("L28").Value2 = 5.479452054794

# ---- Row 31: Hate Crimes ----
$ws.Range("F31").Value2 = 1
